$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9
$ws.Range("D9").Value = "학위 인증 (Accreditation) 후기 – 1"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/accreditation-procedure-experience-1/#utm_source=rss&utm_medium=rss&utm_campaign=accreditation-procedure-experience-1"

# Row 12
$ws.Range("D12").Value = "사이킷런 1.1 버전이 릴리즈되었습니다!"
$ws.Range("E12").Value = "https://tensorflow.blog/2022/05/13/%ec%82%ac%ec%9d%b4%ed%82%b7%eb%9f%b0-1-1-%eb%b2%84%ec%a0%84%ec%9d%b4-%eb%a6%b4%eb%a6%ac%ec%a6%88%eb%90%98%ec%97%88%ec%8a%b5%eb%8b%88%eb%8b%a4/"

# Row 26
$ws.Range("D26").Value = "생성 모델의 새로운 흐름 확산 모델(Diffusion model)에 관하여"

# Row 27
$ws.Range("D27").Value = "개발자를 위한 AWS 클라우드 보안 (1) - 클라우드 설계 원칙과 IAM"
$ws.Range("E27").Value = "https://blog.pingpong.us/aws-cloud-security-for-devs-1/"

# Row 36
$ws.Range("D36").Value = "Various Normalization Techniques for Deep Learning"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/364"

# Row 46
$ws.Range("D46").Value = "패트릭 검사 (Patrick test)"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/462"
